$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.145.94"
$ws.Range("E2").Value = "  -2.03%  "
$ws.Range("D3").Value = "2.437.70"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("E4").Value = "  +0.06%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.44"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  -2.47%  "
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.60"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  -2.17%  "
$ws.Range("E7").Value = "  +0.04%  "
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.498"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  -2.67%  "
$ws.Range("D9").Value = "2.437.20"
$ws.Range("E9").Value = "  -1.87%  "
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.147"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  -6.81%  "
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.163"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  -1.97%  "
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.333"
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = "  -5.92%  "
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.72"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  -3.77%  "
$ws.Range("D14").Value = "2.893.08"
$ws.Range("E14").Value = "  -1.56%  "
$ws.Range("D15").Value = "68.092.50"
$ws.Range("E15").Value = "  -2.02%  "
$ws.Range("E16").Value = "  -5.00%  "
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "22.95"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  -5.40%  "
$ws.Range("D18").Value = "2.444.16"
$ws.Range("E18").Value = "  -1.54%  "
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.75"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  -3.39%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.04"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  -4.34%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "336.18"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  -2.62%  "
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.70"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  -3.86%  "
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  -0.12%  "
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.82"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  -5.15%  "
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.81"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  -4.88%  "
$ws.Range("D26").Value = "2.568.64"
$ws.Range("E26").Value = "  -1.56%  "
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.59"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  -7.50%  "
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  +0.04%  "
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.00"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  -7.33%  "
$ws.Range("D30").Value = "0.0₃0811"
$ws.Range("E30").Value = "  -7.10%  "
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.04"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  -9.04%  "
$ws.Range("E32").Value = "  +0.06%  "
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "420.98"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  -5.19%  "
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.13"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  -5.22%  "
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.62"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  -5.04%  "
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "156.92"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  +0.66%  "
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.98"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("E39").Value = "  -4.58%  "
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.59"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  -3.10%  "
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.299"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  -4.85%  "
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.33"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  -5.27%  "
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.47"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  -6.99%  "
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.06"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.02"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  -6.42%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "132.87"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  -4.68%  "
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.29"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  -4.40%  "
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0710"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  -2.30%  "
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.472"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  -7.84%  "
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.555"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  -3.02%  "
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0899"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  -2.20%  "
